# DataProvider sheet path problem resolved
#
# - Drop the now-unused Sheet2/Sheet3 placeholder tabs (the data provider
#   only ever reads the first sheet).
# - Rename the remaining sheet to "productData" so it matches the path the
#   DataProvider code expects.
# - Move the active selection to C3 (where the next free data row starts).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the extra, empty worksheets.
$null = $wb.Worksheets.Item("Sheet3").Delete()
$null = $wb.Worksheets.Item("Sheet2").Delete()

# Rename the surviving sheet to the name the DataProvider looks for.
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "productData"

# Update the stored selection on the sheet.
$null = $ws.Range("C3").Select()
